# Auto-generated edit script: applies cached-value updates from the commit diff
# to Sheets/Lich_Profits.xlsx (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 363.82352
$ws.Range("I9").Value = 352.46667
$ws.Range("K9").Value = 352.46667
$ws.Range("M9").Value = -183.46667
$ws.Range("H17").Value = 2684085
$ws.Range("J17").Value = 2743729
$ws.Range("L17").Value = 8231187
$ws.Range("N17").Value = -8231523
$ws.Range("H18").Value = 1955.65
$ws.Range("I18").Value = 506.3889
$ws.Range("J18").Value = 14999
$ws.Range("K18").Value = 506.3889
$ws.Range("L18").Value = 14999
$ws.Range("M18").Value = -222.3889
$ws.Range("N18").Value = -15567
$ws.Range("H62").Value = 6416361
$ws.Range("J62").Value = 4778.4
$ws.Range("L62").Value = 4778.4
$ws.Range("N62").Value = -6026.4
$ws.Range("H65").Value = 6416361
$ws.Range("J65").Value = 4778.4
$ws.Range("L65").Value = 23892
$ws.Range("N65").Value = -30132
$ws.Range("H132").Value = 2663
$ws.Range("I132").Value = 2344.3635
$ws.Range("K132").Value = 7033.0905
$ws.Range("M132").Value = -4503.0905
$ws.Range("H138").Value = 3042
$ws.Range("J138").Value = 2935.5227
$ws.Range("L138").Value = 8806.5681
$ws.Range("N138").Value = -19086.5681

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3672.4614
$ws.Range("I28").Value = 3672.4614
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3672.4614
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -3480.4614
$ws.Range("N28").ClearContents()
$ws.Range("H32").Value = 2265.4175
$ws.Range("I32").Value = 1889.6704
$ws.Range("K32").Value = 1889.6704
$ws.Range("M32").Value = -1602.6704
$ws.Range("H99").Value = 3672.4614
$ws.Range("I99").Value = 3672.4614
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3672.4614
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -677.4614000000001
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 3513.5715
$ws.Range("I122").Value = 3398.9312
$ws.Range("K122").Value = 10196.7936
$ws.Range("M122").Value = -7746.793600000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3266.8333
$ws.Range("I20").Value = 3401.0588
$ws.Range("J20").Value = 985
$ws.Range("K20").Value = 3401.0588
$ws.Range("L20").Value = 985
$ws.Range("M20").Value = -3154.0588
$ws.Range("N20").Value = -1479
$ws.Range("H80").Value = 849.12
$ws.Range("I80").Value = 1327.5
$ws.Range("J80").Value = 407.53845
$ws.Range("K80").Value = 1327.5
$ws.Range("L80").Value = 407.53845
$ws.Range("M80").Value = -329.5
$ws.Range("N80").Value = -2403.53845
$ws.Range("H83").Value = 849.12
$ws.Range("I83").Value = 1327.5
$ws.Range("J83").Value = 407.53845
$ws.Range("K83").Value = 6637.5
$ws.Range("L83").Value = 2037.69225
$ws.Range("M83").Value = -1645.5
$ws.Range("N83").Value = -12021.69225
$ws.Range("H99").Value = 3666.3076
$ws.Range("I99").Value = 2912.8096
$ws.Range("J99").Value = 4545.3887
$ws.Range("K99").Value = 2912.8096
$ws.Range("L99").Value = 4545.3887
$ws.Range("M99").Value = -1414.8096
$ws.Range("N99").Value = -7541.3887
$ws.Range("H105").Value = 1863.9678
$ws.Range("I105").Value = 1805.6666
$ws.Range("K105").Value = 1805.6666
$ws.Range("M105").Value = -58.66660000000002
$ws.Range("H132").Value = 94999.25
$ws.Range("J132").Value = 94999.25
$ws.Range("L132").Value = 94999.25
$ws.Range("N132").Value = -105119.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2999.3
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2999.3
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2999.3
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -3339.3
$ws.Range("H68").Value = 37166.418
$ws.Range("J68").Value = 37166.418
$ws.Range("L68").Value = 37166.418
$ws.Range("N68").Value = -38664.418
$ws.Range("H71").Value = 37166.418
$ws.Range("J71").Value = 37166.418
$ws.Range("L71").Value = 111499.254
$ws.Range("N71").Value = -118987.254
$ws.Range("H107").Value = 3805.3618
$ws.Range("I107").Value = 935.125
$ws.Range("J107").Value = 5286.7744
$ws.Range("K107").Value = 935.125
$ws.Range("L107").Value = 5286.7744
$ws.Range("M107").Value = 984.875
$ws.Range("N107").Value = -9126.7744
$ws.Range("H132").Value = 5989.171
$ws.Range("I132").Value = 5864.758
$ws.Range("K132").Value = 17594.274
$ws.Range("M132").Value = -15064.274

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 161.45454
$ws.Range("I2").Value = 209.2
$ws.Range("K2").Value = 1255.2
$ws.Range("M2").Value = -1142.2
$ws.Range("H23").Value = 1078.6428
$ws.Range("I23").Value = 494
$ws.Range("J23").Value = 1517.125
$ws.Range("K23").Value = 1482
$ws.Range("L23").Value = 4551.375
$ws.Range("M23").Value = -1247
$ws.Range("N23").Value = -5021.375
$ws.Range("H38").Value = 172
$ws.Range("I38").Value = 435.8
$ws.Range("K38").Value = 1307.4
$ws.Range("M38").Value = -960.4000000000001
$ws.Range("H40").Value = 91.04761999999999
$ws.Range("I40").Value = 96.4375
$ws.Range("J40").Value = 73.8
$ws.Range("K40").Value = 385.75
$ws.Range("L40").Value = 295.2
$ws.Range("M40").Value = -316.75
$ws.Range("N40").Value = -433.2
$ws.Range("H51").Value = 1624.5
$ws.Range("I51").Value = 1624.5
$ws.Range("K51").Value = 4873.5
$ws.Range("M51").Value = -4413.5
$ws.Range("H55").Value = 8271.154
$ws.Range("I55").Value = 3275
$ws.Range("J55").Value = 9179.546
$ws.Range("K55").Value = 9825
$ws.Range("L55").Value = 27538.638
$ws.Range("M55").Value = -9648
$ws.Range("N55").Value = -27892.638
$ws.Range("H68").Value = 2978071.8
$ws.Range("J68").Value = 2111.739
$ws.Range("L68").Value = 6335.217000000001
$ws.Range("N68").Value = -7957.217000000001
$ws.Range("H71").Value = 2978071.8
$ws.Range("J71").Value = 2111.739
$ws.Range("L71").Value = 19005.651
$ws.Range("N71").Value = -27117.651
$ws.Range("H131").Value = 1454.75
$ws.Range("J131").Value = 1686.5
$ws.Range("L131").Value = 5059.5
$ws.Range("N131").Value = -15139.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 817323.5600000001
$ws.Range("J11").Value = 78666.46000000001
$ws.Range("L11").Value = 78666.46000000001
$ws.Range("N11").Value = -78944.46000000001
$ws.Range("H70").Value = 90918840
$ws.Range("I70").Value = 7995.1665
$ws.Range("K70").Value = 7995.1665
$ws.Range("M70").Value = -7725.1665
$ws.Range("H73").Value = 90918840
$ws.Range("I73").Value = 7995.1665
$ws.Range("K73").Value = 7995.1665
$ws.Range("M73").Value = -7059.1665
$ws.Range("H122").Value = 2818.5715
$ws.Range("I122").Value = 2777.647
$ws.Range("K122").Value = 8332.940999999999
$ws.Range("M122").Value = -5882.940999999999
$ws.Range("H132").Value = 55580.58
$ws.Range("I132").Value = 69389.266
$ws.Range("J132").Value = 3798
$ws.Range("K132").Value = 208167.798
$ws.Range("L132").Value = 11394
$ws.Range("M132").Value = -205637.798
$ws.Range("N132").Value = -16454
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H22").Value = 849.7778
$ws.Range("I22").Value = 774.6667
$ws.Range("K22").Value = 774.6667
$ws.Range("M22").Value = -479.6667
$ws.Range("H27").Value = 849.7778
$ws.Range("I27").Value = 774.6667
$ws.Range("K27").Value = 774.6667
$ws.Range("M27").Value = -667.6667
$ws.Range("H43").Value = 17067.592
$ws.Range("J43").Value = 17088.295
$ws.Range("L43").Value = 17088.295
$ws.Range("N43").Value = -17474.295
$ws.Range("H61").Value = 2489
$ws.Range("J61").Value = 2482
$ws.Range("L61").Value = 2482
$ws.Range("N61").Value = -2886
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H113").Value = 2489
$ws.Range("J113").Value = 2482
$ws.Range("L113").Value = 2482
$ws.Range("N113").Value = -6822
$ws.Range("H122").Value = 4122.385
$ws.Range("J122").Value = 2895
$ws.Range("L122").Value = 8685
$ws.Range("N122").Value = -13585
$ws.Range("H136").Value = 3568.182
$ws.Range("I136").Value = 3528.4348
$ws.Range("J136").Value = 3659.6
$ws.Range("K136").Value = 10585.3044
$ws.Range("L136").Value = 10978.8
$ws.Range("M136").Value = -8035.304400000001
$ws.Range("N136").Value = -16078.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2993
$ws.Range("I126").Value = 2628.3845
$ws.Range("K126").Value = 7885.1535
$ws.Range("M126").Value = -5415.1535
$ws.Range("H132").Value = 20835316
$ws.Range("I132").Value = 34484180
$ws.Range("K132").Value = 103452540
$ws.Range("M132").Value = -103450010
$ws.Range("H136").Value = 347009.3
$ws.Range("I136").Value = 372641.47
$ws.Range("K136").Value = 1117924.41
$ws.Range("M136").Value = -1115374.41
